$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.020.01'
$ws.Range("E2").Value = '  -0.62%  '
$ws.Range("D3").Value = '3.132.26'
$ws.Range("E3").Value = '  -1.49%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '''568.75'
$ws.Range("E5").Value = '  -0.21%  '
$ws.Range("D6").Value = '''161.06'
$ws.Range("E6").Value = '  -5.01%  '
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("E8").Value = '  -7.97%  '
$ws.Range("D9").Value = '''0.116'
$ws.Range("E9").Value = '  -4.29%  '
$ws.Range("D10").Value = '''6.56'
$ws.Range("E10").Value = '  -2.57%  '
$ws.Range("E11").Value = '  -1.99%  '
$ws.Range("D12").Value = '3.681.76'
$ws.Range("E12").Value = '  -1.19%  '
$ws.Range("E13").Value = '  -1.01%  '
$ws.Range("D14").Value = '64.128.97'
$ws.Range("E14").Value = '  -0.51%  '
$ws.Range("D15").Value = '''24.78'
$ws.Range("E15").Value = '  -2.69%  '
$ws.Range("D16").Value = '3.134.25'
$ws.Range("E16").Value = '  -1.49%  '
$ws.Range("E17").Value = '  -3.31%  '
$ws.Range("D18").Value = '''400.37'
$ws.Range("E18").Value = '  -4.75%  '
$ws.Range("D19").Value = '''12.48'
$ws.Range("E19").Value = '  -2.69%  '
$ws.Range("D20").Value = '''5.20'
$ws.Range("E20").Value = '  -2.99%  '
$ws.Range("D21").Value = '''7.08'
$ws.Range("E21").Value = '  +0.49%  '
$ws.Range("E22").Value = '  +3.42%  '
$ws.Range("D23").Value = '''0.999'
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("D24").Value = '''67.99'
$ws.Range("E24").Value = '  -3.34%  '
$ws.Range("D25").Value = '''0.480'
$ws.Range("D26").Value = '''0.193'
$ws.Range("E26").Value = '  -4.89%  '
$ws.Range("D27").Value = '''0.0000100'
$ws.Range("E27").Value = '  -5.60%  '
$ws.Range("D28").Value = '''8.74'
$ws.Range("E28").Value = '  -0.07%  '
$ws.Range("D29").Value = '''0.997'
$ws.Range("E29").Value = '  -0.33%  '
$ws.Range("D30").Value = '''1.80'
$ws.Range("E30").Value = '  -1.65%  '
$ws.Range("D31").Value = '''21.04'
$ws.Range("E31").Value = '  -3.71%  '
$ws.Range("D32").Value = '''6.21'
$ws.Range("D33").Value = '''157.96'
$ws.Range("E33").Value = '  +0.67%  '
$ws.Range("D34").Value = '''4.78'
$ws.Range("E34").Value = '  -4.67%  '
$ws.Range("E35").Value = '  -3.36%  '
$ws.Range("E36").Value = '  -3.33%  '
$ws.Range("D37").Value = '2.666.52'
$ws.Range("E37").Value = '  -1.24%  '
$ws.Range("E38").Value = '  -2.50%  '
$ws.Range("D39").Value = '''23.43'
$ws.Range("E39").Value = '  -4.52%  '
$ws.Range("D40").Value = '''4.05'
$ws.Range("E40").Value = '  -2.79%  '
$ws.Range("D41").Value = '''0.687'
$ws.Range("E41").Value = '  -3.43%  '
$ws.Range("D42").Value = '''0.0609'
$ws.Range("E42").Value = '  -2.05%  '
$ws.Range("D43").Value = '''5.44'
$ws.Range("E43").Value = '  -4.91%  '
$ws.Range("D44").Value = '''0.0254'
$ws.Range("E44").Value = '  -3.05%  '
$ws.Range("D45").Value = '''284.87'
$ws.Range("E45").Value = '  -3.13%  '
$ws.Range("D46").Value = '''20.97'
$ws.Range("E46").Value = '  -3.54%  '
$ws.Range("E47").Value = '  -0.28%  '
$ws.Range("D48").Value = '''0.0971'
$ws.Range("E48").Value = '  -2.19%  '
$ws.Range("D49").Value = '''10.48'
$ws.Range("E49").Value = '  +0.50%  '
$ws.Range("E50").Value = '  -7.02%  '
$ws.Range("D51").Value = '''5.62'
$ws.Range("E51").Value = '  -2.63%  '
